$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111354726
$ws.Range("B2").Value = 96348
$ws.Range("D2").Value = 'VU'
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = 'Knärot'
$ws.Range("G2").Value = 'Goodyera repens'
$ws.Range("H2").Value = '(L.) R. Br.'
$ws.Range("Q2").Value = 555138.3841980774
$ws.Range("R2").Value = 6999122.632572887
$ws.Range("Z2").Value = '18:59'
$ws.Range("AB2").Value = '18:59'
$ws.Range("AC2").Value = 'Rikligt'

# Row 3
$ws.Range("A3").Value = 111354546
$ws.Range("B3").Value = 77515
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = 'Garnlav'
$ws.Range("G3").Value = 'Alectoria sarmentosa'
$ws.Range("H3").Value = '(Ach.) Ach.'
$ws.Range("Q3").Value = 555144.7642560177
$ws.Range("R3").Value = 6999173.89047078

# Row 4
$ws.Range("A4").Value = 111354025
$ws.Range("B4").Value = 89405
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = 'Ullticka'
$ws.Range("G4").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H4").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q4").Value = 555154.7369911602
$ws.Range("R4").Value = 6999253.724818715
$ws.Range("Z4").Value = '00:00'
$ws.Range("AB4").Value = '00:00'

# Row 5
$ws.Range("A5").Value = 111355331
$ws.Range("Q5").Value = 555224.1634512447
$ws.Range("R5").Value = 6999062.984709017

# Row 6
$ws.Range("A6").Value = 111355391
$ws.Range("Q6").Value = 555215.2391852245
$ws.Range("R6").Value = 6999054.684715276

# Row 7
$ws.Range("A7").Value = 111354809
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("Q7").Value = 555164.2519777509
$ws.Range("R7").Value = 6999120.809599082
$ws.Range("Z7").Value = '19:06'
$ws.Range("AB7").Value = '19:06'
$ws.Range("AC7").Value = 'Rikligt'

# Row 8
$ws.Range("A8").Value = 111355282
$ws.Range("Q8").Value = 555252.2533465028
$ws.Range("R8").Value = 6999037.209103072

# Row 10
$ws.Range("A10").Value = 111355197
$ws.Range("B10").Value = 78578
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 6458
$ws.Range("F10").Value = 'Lunglav'
$ws.Range("G10").Value = 'Lobaria pulmonaria'
$ws.Range("H10").Value = '(L.) Hoffm.'
$ws.Range("Q10").Value = 555280.6270040129
$ws.Range("R10").Value = 6999021.397055306
$ws.Range("AC10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111355254
$ws.Range("B11").Value = 78578
$ws.Range("E11").Value = 6458
$ws.Range("F11").Value = 'Lunglav'
$ws.Range("G11").Value = 'Lobaria pulmonaria'
$ws.Range("H11").Value = '(L.) Hoffm.'
$ws.Range("Q11").Value = 555272.5859257083
$ws.Range("R11").Value = 6999014.46986856
$ws.Range("Z11").Value = '19:06'
$ws.Range("AB11").Value = '19:06'

# Row 12
$ws.Range("A12").Value = 111355227
$ws.Range("B12").Value = 78578
$ws.Range("E12").Value = 6458
$ws.Range("F12").Value = 'Lunglav'
$ws.Range("G12").Value = 'Lobaria pulmonaria'
$ws.Range("H12").Value = '(L.) Hoffm.'
$ws.Range("Q12").Value = 555268.7908350837
$ws.Range("R12").Value = 6999024.363821026
$ws.Range("Z12").Value = '19:06'
$ws.Range("AB12").Value = '19:06'
$ws.Range("AC12").ClearContents()

# Row 13
$ws.Range("A13").Value = 111354189
$ws.Range("Q13").Value = 555134.7763198819
$ws.Range("R13").Value = 6999254.742998262
$ws.Range("Z13").Value = '00:00'
$ws.Range("AB13").Value = '00:00'
